$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Worksheet data: replace numeric index column A with month labels,
#     and update column B values; extend the table from 7 to 12 rows ---
$months = @(
  "December 2015", "January 2016", "February 2016", "March 2016",
  "April 2016", "May 2016", "June 2016", "July 2016",
  "August 2016", "September 2016", "October 2016", "November 2016"
)
$values = @(578, 527, 424, 427, 538, 557, 488, 454, 651, 488, 426, 448)

for ($i = 0; $i -lt $months.Count; $i++) {
  $row = $i + 1

  # Force text entry for column A so "December 2015" etc. aren't
  # auto-converted to date serials, then strip the temporary number
  # format so no style residue is left on the cell.
  $cellA = $ws.Cells.Item($row, 1)
  $cellA.NumberFormat = "@"
  $cellA.Value = $months[$i]
  $cellA.ClearFormats()

  $ws.Cells.Item($row, 2).Value = $values[$i]
}

# --- Chart updates ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart

# New chart style (xl/charts/chart1.xml: <style val="4"/> -> <style val="8"/>)
$chart.ChartStyle = 8

# Manual plot-area layout sized to 85% (factor mode) in both dimensions
$pa = $chart.PlotArea
$pa.Width = 0.85
$pa.Height = 0.85

# Extend the series' category/value references from row 7 to row 12
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(,'Sheet1'!`$A`$1:`$A`$12,'Sheet1'!`$B`$1:`$B`$12,1)"

Write-Host "done"
